# Auto-generated edit script to update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.403.56"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.574.00"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'291.28"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.3767"
$ws.Range("E7").Value = "  +2.84%  "
$ws.Range("D8").Value = "'49.89"
$ws.Range("D9").Value = "'0.3413"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").Value = "'1.165"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").Value = "'0.07672"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "'21.35"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "'5.982"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "'6.920"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "1.571.78"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "'90.58"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "'0.06741"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'16.75"
$ws.Range("E21").Value = "  +2.57%  "
$ws.Range("D22").Value = "'6.229"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'0.5278"
$ws.Range("E23").Value = "  -4.28%  "
$ws.Range("D24").Value = "'12.01"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").Value = "22.410.66"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "'2.424"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").Value = "'2.769"
$ws.Range("E27").Value = "  -7.01%  "
$ws.Range("D28").Value = "'20.28"
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("D29").Value = "'145.07"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "'5.056"
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("D31").Value = "'126.18"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").Value = "1.749.04"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "'6.213"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'1.016"
$ws.Range("E34").Value = "  +4.70%  "
$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'2.021"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("D36").Value = "'10.06"
$ws.Range("E36").Value = "  -3.78%  "
$ws.Range("D37").Value = "'0.08561"
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("D38").Value = "'0.02567"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").Value = "'1.340"
$ws.Range("E40").Value = "  +6.06%  "
$ws.Range("D41").Value = "'0.06525"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'5.480"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.6489"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'11.63"
$ws.Range("E44").Value = "  -1.90%  "
$ws.Range("D45").Value = "'14.10"
$ws.Range("E45").Value = "  -4.01%  "
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").Value = "'0.6042"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D48").Value = "'3.787"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "'1.301"
$ws.Range("E49").Value = "  +9.79%  "
$ws.Range("D50").Value = "'2.099"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").Value = "'125.33"
$ws.Range("E51").Value = "  +3.01%  "
